$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(2, 3).Value = 316056
$ws.Cells.Item(2, 4).Value = 402889751
$ws.Cells.Item(3, 3).Value = 255
$ws.Cells.Item(3, 4).Value = 303979
$ws.Cells.Item(8, 3).Value = 850
$ws.Cells.Item(8, 4).Value = 1250408
$ws.Cells.Item(10, 3).Value = 116031
$ws.Cells.Item(10, 4).Value = 170024789
$ws.Cells.Item(12, 3).Value = 58567
$ws.Cells.Item(12, 4).Value = 84530270
$ws.Cells.Item(16, 3).Value = 3969
$ws.Cells.Item(16, 4).Value = 5633019
$ws.Cells.Item(20, 3).Value = 6477
$ws.Cells.Item(20, 4).Value = 9037598
$ws.Cells.Item(22, 3).Value = 76420
$ws.Cells.Item(22, 4).Value = 95365815
$ws.Cells.Item(28, 3).Value = 32200
$ws.Cells.Item(28, 4).Value = 47141858
$ws.Cells.Item(30, 3).Value = 11355
$ws.Cells.Item(30, 4).Value = 16332147
$ws.Cells.Item(35, 3).Value = 1779
$ws.Cells.Item(35, 4).Value = 2510333
$ws.Cells.Item(36, 3).Value = 95987
$ws.Cells.Item(36, 4).Value = 120901635
$ws.Cells.Item(37, 3).Value = 66
$ws.Cells.Item(37, 4).Value = 77037
$ws.Cells.Item(44, 3).Value = 44048
$ws.Cells.Item(44, 4).Value = 64557830
$ws.Cells.Item(45, 3).Value = 26
$ws.Cells.Item(45, 4).Value = 38950
$ws.Cells.Item(46, 3).Value = 9027
$ws.Cells.Item(46, 4).Value = 12955698
$ws.Cells.Item(51, 3).Value = 2253
$ws.Cells.Item(51, 4).Value = 3143001
$ws.Cells.Item(52, 3).Value = 68135
$ws.Cells.Item(52, 4).Value = 85513101
$ws.Cells.Item(56, 3).Value = 381
$ws.Cells.Item(56, 4).Value = 559965
$ws.Cells.Item(58, 3).Value = 27903
$ws.Cells.Item(58, 4).Value = 40922626
$ws.Cells.Item(61, 3).Value = 10954
$ws.Cells.Item(61, 4).Value = 15837956
$ws.Cells.Item(67, 3).Value = 1439
$ws.Cells.Item(67, 4).Value = 2014921
$ws.Cells.Item(69, 3).Value = 20239
$ws.Cells.Item(69, 4).Value = 26511372
$ws.Cells.Item(73, 3).Value = 7512
$ws.Cells.Item(73, 4).Value = 10997530
$ws.Cells.Item(75, 3).Value = 5047
$ws.Cells.Item(75, 4).Value = 7328606
$ws.Cells.Item(77, 3).Value = 268
$ws.Cells.Item(77, 4).Value = 375673
$ws.Cells.Item(78, 3).Value = 138829
$ws.Cells.Item(78, 4).Value = 173173406
$ws.Cells.Item(84, 3).Value = 62976
$ws.Cells.Item(84, 4).Value = 92306039
$ws.Cells.Item(87, 3).Value = 29339
$ws.Cells.Item(87, 4).Value = 42442007
$ws.Cells.Item(89, 3).Value = 2714
$ws.Cells.Item(89, 4).Value = 3908360
$ws.Cells.Item(90, 3).Value = 2762
$ws.Cells.Item(90, 4).Value = 3903280
$ws.Cells.Item(91, 3).Value = 32054
$ws.Cells.Item(91, 4).Value = 43432262
$ws.Cells.Item(95, 3).Value = 7792
$ws.Cells.Item(95, 4).Value = 11457470
$ws.Cells.Item(97, 3).Value = 7088
$ws.Cells.Item(97, 4).Value = 10275127
$ws.Cells.Item(101, 3).Value = 8791
$ws.Cells.Item(101, 4).Value = 12200726
$ws.Cells.Item(103, 3).Value = 2210
$ws.Cells.Item(103, 4).Value = 3256352
$ws.Cells.Item(105, 3).Value = 2979
$ws.Cells.Item(105, 4).Value = 4351241
$ws.Cells.Item(107, 3).Value = 128
$ws.Cells.Item(107, 4).Value = 186120
$ws.Cells.Item(108, 3).Value = 174
$ws.Cells.Item(108, 4).Value = 247086
$ws.Cells.Item(109, 3).Value = 139319
$ws.Cells.Item(109, 4).Value = 172303997
$ws.Cells.Item(113, 3).Value = 948
$ws.Cells.Item(113, 4).Value = 1392288
$ws.Cells.Item(115, 3).Value = 52224
$ws.Cells.Item(115, 4).Value = 76559046
$ws.Cells.Item(117, 3).Value = 26595
$ws.Cells.Item(117, 4).Value = 38531293
$ws.Cells.Item(118, 3).Value = 1300
$ws.Cells.Item(118, 4).Value = 1779551
$ws.Cells.Item(121, 3).Value = 2201
$ws.Cells.Item(121, 4).Value = 3090850
$ws.Cells.Item(123, 3).Value = 493662
$ws.Cells.Item(123, 4).Value = 651039269
$ws.Cells.Item(124, 3).Value = 89
$ws.Cells.Item(124, 4).Value = 117789
$ws.Cells.Item(128, 3).Value = 1363
$ws.Cells.Item(128, 4).Value = 2020811
$ws.Cells.Item(130, 3).Value = 204798
$ws.Cells.Item(130, 4).Value = 301066008
$ws.Cells.Item(133, 3).Value = 176918
$ws.Cells.Item(133, 4).Value = 257153645
$ws.Cells.Item(136, 3).Value = 2809
$ws.Cells.Item(136, 4).Value = 3948685
$ws.Cells.Item(138, 3).Value = 6168
$ws.Cells.Item(138, 4).Value = 8715072
$ws.Cells.Item(141, 3).Value = 43790
$ws.Cells.Item(141, 4).Value = 58473348
$ws.Cells.Item(147, 3).Value = 13896
$ws.Cells.Item(147, 4).Value = 20382249
$ws.Cells.Item(148, 3).Value = 3696
$ws.Cells.Item(148, 4).Value = 5330140
$ws.Cells.Item(154, 3).Value = 17254
$ws.Cells.Item(154, 4).Value = 22797660
$ws.Cells.Item(158, 3).Value = 7046
$ws.Cells.Item(158, 4).Value = 10247481
$ws.Cells.Item(160, 3).Value = 4909
$ws.Cells.Item(160, 4).Value = 7066936
$ws.Cells.Item(165, 3).Value = 15353
$ws.Cells.Item(165, 4).Value = 22279946
$ws.Cells.Item(166, 3).Value = 1737
$ws.Cells.Item(166, 4).Value = 2584030
$ws.Cells.Item(167, 3).Value = 234
$ws.Cells.Item(167, 4).Value = 345802
$ws.Cells.Item(169, 3).Value = 52
$ws.Cells.Item(169, 4).Value = 77690
$ws.Cells.Item(171, 3).Value = 86557
$ws.Cells.Item(171, 4).Value = 108283751
$ws.Cells.Item(178, 3).Value = 33535
$ws.Cells.Item(178, 4).Value = 49181447
$ws.Cells.Item(180, 3).Value = 12831
$ws.Cells.Item(180, 4).Value = 18537538
$ws.Cells.Item(186, 3).Value = 235153
$ws.Cells.Item(186, 4).Value = 292370204
$ws.Cells.Item(188, 3).Value = 166
$ws.Cells.Item(188, 4).Value = 239236
$ws.Cells.Item(194, 3).Value = 85860
$ws.Cells.Item(194, 4).Value = 125866472
$ws.Cells.Item(197, 3).Value = 32619
$ws.Cells.Item(197, 4).Value = 46945523
$ws.Cells.Item(200, 3).Value = 5048
$ws.Cells.Item(200, 4).Value = 7193913
$ws.Cells.Item(203, 3).Value = 4738
$ws.Cells.Item(203, 4).Value = 6556383
$ws.Cells.Item(206, 3).Value = 260036
$ws.Cells.Item(206, 4).Value = 321877783
$ws.Cells.Item(215, 3).Value = 94262
$ws.Cells.Item(215, 4).Value = 137905682
$ws.Cells.Item(218, 3).Value = 50761
$ws.Cells.Item(218, 4).Value = 73361684
$ws.Cells.Item(221, 3).Value = 4626
$ws.Cells.Item(221, 4).Value = 6493750
$ws.Cells.Item(224, 3).Value = 5588
$ws.Cells.Item(224, 4).Value = 7731405
$ws.Cells.Item(227, 3).Value = 104766
$ws.Cells.Item(227, 4).Value = 131118755
$ws.Cells.Item(234, 3).Value = 49053
$ws.Cells.Item(234, 4).Value = 71866331
$ws.Cells.Item(236, 3).Value = 12216
$ws.Cells.Item(236, 4).Value = 17562969
$ws.Cells.Item(238, 3).Value = 1881
$ws.Cells.Item(238, 4).Value = 2696338
$ws.Cells.Item(240, 3).Value = 2433
$ws.Cells.Item(240, 4).Value = 3399137
$ws.Cells.Item(241, 3).Value = 253670
$ws.Cells.Item(241, 4).Value = 320355652
$ws.Cells.Item(249, 3).Value = 94820
$ws.Cells.Item(249, 4).Value = 138944337
$ws.Cells.Item(252, 3).Value = 64039
$ws.Cells.Item(252, 4).Value = 92800563
$ws.Cells.Item(254, 3).Value = 2386
$ws.Cells.Item(254, 4).Value = 3366428
$ws.Cells.Item(257, 3).Value = 4498
$ws.Cells.Item(257, 4).Value = 6314164
